$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 64 data: second CRM test run, 11/12/2019 (accounting for
# potentially bad acid influence)

# Date (copy date format from A63 so we reuse the existing style index
# instead of minting a brand-new number format)
$ws.Range("A63").Copy()
$ws.Range("A64").PasteSpecial(-4122)
$ws.Range("A64").Value = 43781

$ws.Range("B64").Value = 2212.43486397574
$ws.Range("C64").Value = 2207.0300000000002

# % off = 100*(B64-C64)/C64, continuing the column's running formula
$ws.Range("D64").Formula = "=100*(B64-C64)/C64"

$ws.Range("E64").Value = 169

$ws.Range("F64").Value = "second test run of the day"

# Reflect the sheet scrolling down one row and the new active selection
$ws.Range("F65").Select() | Out-Null
